# Mar 25 - Input Update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update shared string text in place: 'deluxe25off' -> 'deluxe25offp-redes'
# (handled implicitly by writing the new text value into the cell that used it)

# Build full data array for rows 2-41, columns A-E
$data = @(
    @("QA","PrincipalSecret","Core","Kit","Chrome"),
    @("QA","TryDermaFlash","Core","Kit","Chrome"),
    @("QA","Sub-D","Core","Kit","Chrome"),
    @("QA","SheerCover","Core","Kit","Chrome"),
    @("QA","SpecificBeauty","Core","Kit","Chrome"),
    @("QA","PrincipalSecret","Order30","Kit","Chrome"),
    @("QA","Sub-D","deluxe25offp","Kit","Chrome"),
    @("QA","ReclaimBotanical","Core","Kit","Chrome"),
    @("QA","SpecificBeauty","deluxe-offer","Kit","Chrome"),
    @("QA","TryDermaFlash","pnln","Kit","Chrome"),
    @("QA","Sub-D","cpcb2017","Kit","Chrome"),
    @("QA","TryDermaFlash","trydermaflash","Kit","Chrome"),
    @("QA","Sub-D","deluxe25offp-redes","Kit","Chrome"),
    @("QA","Sub-D","cpwbunusedbdbj","Kit","Chrome"),
    @("End","","","",""),
    @("","","","",""),
    @("","","","",""),
    @("","","","",""),
    @("","","","",""),
    @("","","","",""),
    @("","","","",""),
    @("","","","",""),
    @("","","","",""),
    @("","","","",""),
    @("Environment","Brand","Campaign","Categories","Browser"),
    @("QA","PrincipalSecret","Core","Kit","Chrome"),
    @("QA","PrincipalSecret","Order30","Kit","Chrome"),
    @("QA","ReclaimBotanical","Core","Kit","Chrome"),
    @("QA","SheerCover","Core","Kit","Chrome"),
    @("QA","SpecificBeauty","Core","Kit","Chrome"),
    @("QA","SpecificBeauty","deluxe-offer","Kit","Chrome"),
    @("QA","Sub-D","Core","Kit","Chrome"),
    @("QA","Sub-D","cpcb2017","Kit","Chrome"),
    @("QA","Sub-D","deluxe25offp-redes","Kit","Chrome"),
    @("QA","Sub-D","deluxe25offp","Kit","Chrome"),
    @("QA","Sub-D","cpwbunusedbdbj","Kit","Chrome"),
    @("QA","TryDermaFlash","Core","Kit","Chrome"),
    @("QA","TryDermaFlash","pnln","Kit","Chrome"),
    @("QA","TryDermaFlash","trydermaflash","Kit","Chrome"),
    @("End","","","","")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $val = $rowVals[$j]
        if ($val -eq "") {
            $ws.Cells.Item($rowNum, $colNum).Value = $null
        } else {
            $ws.Cells.Item($rowNum, $colNum).Value = $val
        }
    }
}

Write-Output "Applied data update"